$d = $word.ActiveDocument

# Several paragraphs in this document were corrupted by a bug that
# injected a stray table (<w:tbl>) directly inside a run (<w:r>) - an
# invalid place for a table to live in OOXML. Each such paragraph is
# wrapped in a <w:bookmarkStart>/<w:bookmarkEnd> pair whose name
# duplicates a legitimate bookmark that already exists elsewhere in the
# document around real prose. Because the table sits in that invalid
# location, it is invisible to the normal text-oriented object model
# (Range.Text, Find.Execute, Tables, ...): those APIs see these
# paragraphs as empty. They are however still real entries in
# Paragraphs, and Range.WordOpenXML faithfully reveals the embedded
# <w:tbl> markup, which lets us find and remove them precisely.

$tableParagraphIndexes = New-Object System.Collections.ArrayList

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $openXml = $para.Range.WordOpenXML

    # Range.WordOpenXML returns a full mini-package whose <w:body> begins
    # with exactly this paragraph's own markup, followed by an
    # auto-generated trailing paragraph/sectPr. Looking only at the start
    # of <w:body> tells us whether *this* paragraph itself contains the
    # offending table, rather than matching unrelated <w:tbl> references
    # that legitimately show up later in the package (e.g. table styles).
    $bodyIndex = $openXml.IndexOf("<w:body>")
    $isStrayTableParagraph = $false
    if ($bodyIndex -ge 0) {
        $windowLength = [Math]::Min(800, $openXml.Length - $bodyIndex)
        $paragraphStart = $openXml.Substring($bodyIndex, $windowLength)
        if ($paragraphStart -like "*<w:tbl>*") {
            $isStrayTableParagraph = $true
        }
    }

    if ($isStrayTableParagraph) {
        [void]$tableParagraphIndexes.Add($i)
    }
}

# Delete from the highest index down to the lowest so the indexes of
# paragraphs we have not processed yet stay valid while we work.
for ($k = $tableParagraphIndexes.Count - 1; $k -ge 0; $k--) {
    $idx = $tableParagraphIndexes[$k]
    $para = $d.Paragraphs.Item($idx)
    $para.Range.Delete()
}
